# Slide 11, "Content Placeholder 2" shape, 3rd paragraph:
#   "lw/sw memory addresses are computed beforehand and accessed from the
#    computation register during memory access instructions"
# becomes:
#   "Load word/store word memory addresses are computed beforehand and
#    accessed from the computation register during memory access instructions"
#
# The paragraph originally has 7 runs:
#   1: "l"                                              (dirty=0, err=1)
#   2: "w"                                               (dirty=0, err=1, smtClean=0)
#   3: "/"                                                (dirty=0, smtClean=0)
#   4: "sw"                                               (dirty=0, err=1, smtClean=0)
#   5: " memory addresses are computed beforehand "       (dirty=0, smtClean=0)
#   6: "and accessed "                                    (smtClean=0)
#   7: "from the computation register during memory access instructions" (dirty=0, smtClean=0)
#
# We keep runs 3, 6 and 7 (whose existing formatting already matches what the
# final three runs need) and change their text, then clear (which removes)
# runs 1, 2, 4 and 5.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)

# Re-text the three runs we are keeping first, while the indices still refer
# to the original 7-run layout.
$para.Runs(7, 1).Text = "memory addresses are computed beforehand and accessed from the computation register during memory access instructions"
$para.Runs(6, 1).Text = "word/store word "
$para.Runs(3, 1).Text = "Load "

# Clear out the runs we no longer need. Clearing a run's text removes it from
# the run collection, so we work from the highest index down to keep the
# remaining indices stable.
$para.Runs(5, 1).Text = ""
$para.Runs(4, 1).Text = ""
$para.Runs(2, 1).Text = ""
$para.Runs(1, 1).Text = ""
